$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to D and E columns (rows 2-51) to guard against
# Excel auto-converting numeric-looking / date-looking strings.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '24.033.49'
$ws.Cells.Item(2, 5).Value = '  -4.18%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.638.45'
$ws.Cells.Item(3, 5).Value = '  -3.94%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '1.003'
$ws.Cells.Item(4, 5).Value = '  -0.03%  '

# Row 5
$ws.Cells.Item(5, 2).Value = 'USDC'
$ws.Cells.Item(5, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Cells.Item(5, 4).Value = '1.003'
$ws.Cells.Item(5, 5).Value = '  +0.05%  '

# Row 6
$ws.Cells.Item(6, 2).Value = 'BNB'
$ws.Cells.Item(6, 3).Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Cells.Item(6, 4).Value = '306.68'
$ws.Cells.Item(6, 5).Value = '  -3.11%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '0.3896'
$ws.Cells.Item(7, 5).Value = '  -2.72%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '0.3820'
$ws.Cells.Item(8, 5).Value = '  -5.40%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '1.004'
$ws.Cells.Item(9, 5).Value = '  -0.08%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '49.21'
$ws.Cells.Item(10, 5).Value = '  -6.96%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '1.332'
$ws.Cells.Item(11, 5).Value = '  -9.40%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '0.08367'
$ws.Cells.Item(12, 5).Value = '  -5.30%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '23.50'
$ws.Cells.Item(13, 5).Value = '  -9.45%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '7.038'
$ws.Cells.Item(14, 5).Value = '  -6.09%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '0.00001267'
$ws.Cells.Item(15, 5).Value = '  -6.55%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '7.410'
$ws.Cells.Item(16, 5).Value = '  -6.88%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '1.630.50'
$ws.Cells.Item(17, 5).Value = '  -4.42%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '94.61'
$ws.Cells.Item(18, 5).Value = '  -1.79%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '0.06863'
$ws.Cells.Item(19, 5).Value = '  -4.65%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '21.01'
$ws.Cells.Item(20, 5).Value = '  +1.14%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '6.862'
$ws.Cells.Item(21, 5).Value = '  -5.63%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '1.002'
$ws.Cells.Item(22, 5).Value = '  -0.07%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '13.48'
$ws.Cells.Item(23, 5).Value = '  -5.99%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '24.027.19'
$ws.Cells.Item(24, 5).Value = '  -4.20%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -2.99%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '2.672'
$ws.Cells.Item(26, 5).Value = '  -9.04%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '22.16'
$ws.Cells.Item(27, 5).Value = '  -6.13%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '156.62'
$ws.Cells.Item(28, 5).Value = '  -3.83%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '8.594'
$ws.Cells.Item(29, 5).Value = '  +3.01%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '139.58'
$ws.Cells.Item(30, 5).Value = '  -7.89%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '5.295'
$ws.Cells.Item(31, 5).Value = '  -14.79%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '2.429'
$ws.Cells.Item(32, 5).Value = '  -8.97%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '1.823.23'
$ws.Cells.Item(33, 5).Value = '  -3.78%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '6.809'
$ws.Cells.Item(34, 5).Value = '  -5.32%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '0.07919'
$ws.Cells.Item(35, 5).Value = '  -7.85%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '0.02872'
$ws.Cells.Item(36, 5).Value = '  -9.11%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'Algorand'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(37, 4).Value = '0.2657'
$ws.Cells.Item(37, 5).Value = '  -8.79%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'ImmutableX'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(38, 4).Value = '0.9387'
$ws.Cells.Item(38, 5).Value = '  -10.18%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.09140'
$ws.Cells.Item(39, 5).Value = '  -6.32%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '1.438'
$ws.Cells.Item(40, 5).Value = '  -2.61%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '9.823'
$ws.Cells.Item(41, 5).Value = '  -10.63%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '0.7470'
$ws.Cells.Item(42, 5).Value = '  -9.93%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '12.88'
$ws.Cells.Item(43, 5).Value = '  -8.37%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '15.83'
$ws.Cells.Item(44, 5).Value = '  -7.40%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '0.6819'
$ws.Cells.Item(45, 5).Value = '  -8.12%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '2.441'
$ws.Cells.Item(46, 5).Value = '  -8.84%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '4.080'
$ws.Cells.Item(47, 5).Value = '  -4.05%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  -0.10%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '0.08316'
$ws.Cells.Item(49, 5).Value = '  -8.63%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '131.23'
$ws.Cells.Item(50, 5).Value = '  -6.46%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '1.238'
$ws.Cells.Item(51, 5).Value = '  -13.83%  '
